$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Id" column (column A) - shifts "Level"/"Indicator" left
# by one column and drops the old column C entirely, matching the
# sharedStrings / sheetData / dataValidation shrink in the diff.
$ws.Range("A1").EntireColumn.Delete()

# Post-edit selection moves to C1 (first column past the new data,
# matching the <selection activeCell="C1" sqref="C1"/> in the diff).
$ws.Range("C1").Select()
